# Add a new Job Posting row (Job_Id = JD_001) to the LinkedIn job posting sheet.
# Final layout:
#   Row 1 (headers, bold + bordered "header" style): Job_Id, Job_Title, Job_Description,
#          Total_Years_Min_Exp, Total_Years_Max_Exp, Linked_Poster, Linked_Posted,
#          Resume_received, Resume_downloaded
#   Row 2 (data, default/unstyled): JD_001, Senior Engineer, <job description>, 1, 4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:I1) ---------------------------------------------------
$headers = @(
    "Job_Id",
    "Job_Title",
    "Job_Description",
    "Total_Years_Min_Exp",
    "Total_Years_Max_Exp",
    "Linked_Poster",
    "Linked_Posted",
    "Resume_received",
    "Resume_downloaded"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# A1:C1 already carry the bold/centered/bordered header style from the template;
# copy that formatting across the newly added header cells (D1:I1) so the whole
# header row looks consistent.
$ws.Range("A1").Copy()
$ws.Range("D1:I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data row (A2:E2) ------------------------------------------------------
$ws.Range("A2").Value = "JD_001"
# The template's A2 carried a "wrap text" style; the new data row uses the
# default (unstyled) look, so clear it back to Normal.
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "Senior Engineer"

$jobDescription = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`n" +
    "Work with global teams to drive innovation and deliver scalable applications.`n" +
    "Join Akkodis and be part of a tech-driven, collaborative environment."
$ws.Range("C2").Value = $jobDescription

$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4

# Entering the multi-line description auto-expands row 2's height; restore the
# row to its normal (non-custom) auto height to match the sheet's default look.
$ws.Rows.Item(2).AutoFit()
